$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1
$ws.Range("B1").Value = 2
$ws.Range("C1").Value = 3

$ws.Range("A2").Value = 123
$ws.Range("B2").Value = 123
$ws.Range("C2").Value = 123

$ws.Range("C1:C2").NumberFormat = "0"

$ws.Range("E5").Select()
